$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column J ("Misc. Data") header and blank values for all data rows
$ws.Range("J1").Value = "Misc. Data"
$ws.Range("J2").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("J6").Value = ""
$ws.Range("J7").Value = ""
$ws.Range("J8").Value = ""
$ws.Range("J9").Value = ""
$ws.Range("J10").Value = ""
$ws.Range("J11").Value = ""

# Fix the Authors field (column E) for rows whose value had accumulated
# whitespace padding from repeated failed update attempts; use the final,
# fully-padded author string that the fixed pipeline now successfully applies.
$ws.Range("E2").Value = '[Jianlei%Cao%NULL%0,                          Xiaorong%Hu%NULL%1,                          Wenlin%Cheng%NULL%6,                          Lei%Yu%NULL%6,                          Wen-Jun%Tu%tuwenjun@irm-cams.ac.cn%6,                          Qiang%Liu%liuqiang@irm-cams.ac.cn%7,                          Qiang%Liu%liuqiang@irm-cams.ac.cn%0]'
$ws.Range("E4").Value = '[Chaolin%Huang%NULL%0,                          Yeming%Wang%NULL%14,                          Xingwang%Li%NULL%6,                          Lili%Ren%NULL%5,                          Jianping%Zhao%NULL%12,                          Yi%Hu%NULL%11,                          Li%Zhang%NULL%11,                          Guohui%Fan%NULL%14,                          Jiuyang%Xu%NULL%14,                          Xiaoying%Gu%NULL%14,                          Zhenshun%Cheng%NULL%7,                          Ting%Yu%NULL%24,                          Jiaan%Xia%NULL%6,                          Yuan%Wei%NULL%18,                          Wenjuan%Wu%NULL%5,                          Xuelei%Xie%NULL%6,                          Wen%Yin%NULL%9,                          Hui%Li%NULL%15,                          Min%Liu%NULL%5,                          Yan%Xiao%NULL%5,                          Hong%Gao%NULL%5,                          Li%Guo%NULL%6,                          Jungang%Xie%NULL%6,                          Guangfa%Wang%NULL%5,                          Rongmeng%Jiang%NULL%5,                          Zhancheng%Gao%NULL%5,                          Qi%Jin%NULL%5,                          Jianwei%Wang%wangjw28@163.com%5,                          Bin%Cao%caobin_ben@163.com%13]'
$ws.Range("E6").Value = '[Qiurong%Ruan%NULL%0,                          Kun%Yang%NULL%8,                          Kun%Yang%NULL%0,                          Wenxia%Wang%NULL%8,                          Wenxia%Wang%NULL%0,                          Lingyu%Jiang%NULL%8,                          Lingyu%Jiang%NULL%0,                          Jianxin%Song%songsingsjx@sina.com%8,                          Jianxin%Song%songsingsjx@sina.com%0]'
$ws.Range("E7").Value = '[Carb\u00f3-Cisnero%Yacquelin%coreGivesNoEmail%0,                         Fern\u00e1ndez-Gonz\u00e1lez%Paula%coreGivesNoEmail%3,                         Hierrezuelo-Rojas%Naifi%coreGivesNoEmail%3,                         Subert-Salas%Lizandra%coreGivesNoEmail%3]'
$ws.Range("E8").Value = '[Aguilar%Medina Jos\u00e9 M.%coreGivesNoEmail%0,                         Claro%Vald\u00e9s Ram\u00f3n%coreGivesNoEmail%2,                         Obreg\u00f3n%Santos Angel G.%coreGivesNoEmail%2,                         Proh\u00edas%Mart\u00ednez Juan%coreGivesNoEmail%2,                         Rodr\u00edguez%Blanco Suilbert%coreGivesNoEmail%2]'
$ws.Range("E10").Value = '[Xiaobo%Yang%NULL%0,                          Yuan%Yu%NULL%8,                          Jiqian%Xu%NULL%8,                          Huaqing%Shu%NULL%6,                          Jia''an%Xia%NULL%10,                          Hong%Liu%NULL%6,                          Yongran%Wu%NULL%8,                          Lu%Zhang%NULL%6,                          Zhui%Yu%NULL%6,                          Minghao%Fang%NULL%6,                          Ting%Yu%NULL%0,                          Yaxin%Wang%NULL%8,                          Shangwen%Pan%NULL%6,                          Xiaojing%Zou%NULL%6,                          Shiying%Yuan%NULL%6,                          You%Shang%NULL%8]'
$ws.Range("E11").Value = '[Fei%Zhou%NULL%0,                          Ting%Yu%NULL%0,                          Ronghui%Du%NULL%8,                          Guohui%Fan%NULL%0,                          Ying%Liu%NULL%13,                          Zhibo%Liu%NULL%9,                          Jie%Xiang%NULL%15,                          Yeming%Wang%NULL%0,                          Bin%Song%NULL%9,                          Xiaoying%Gu%NULL%0,                          Lulu%Guan%NULL%8,                          Yuan%Wei%NULL%0,                          Hui%Li%NULL%0,                          Xudong%Wu%NULL%9,                          Jiuyang%Xu%NULL%0,                          Shengjin%Tu%NULL%8,                          Yi%Zhang%NULL%9,                          Hua%Chen%NULL%9,                          Bin%Cao%NULL%0]'
